$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data (weekly update): rows 16..51 of columns D (Fecha), L (Calidad),
# M (Volumen), N (Precio minimo), O (Precio maximo), P (Precio promedio
# ponderado), S (Precio $/Kg).  This reflects: the oldest record (old row 16)
# is dropped, every remaining record shifts up one row, and two brand new
# records are appended (new rows 39-40, i.e. immediately before the old
# tail that already continues the sequence), which ultimately pushes the
# whole block down so the table now spans through row 51 instead of row 50.
$data = @(
    @(16, 44685, "Primera", 160, 15000, 16000, 15500, 861),
    @(17, 44609, "Primera", 240, 13000, 14000, 13500, 750),
    @(18, 44609, "Segunda", 240, 11000, 12000, 11500, 639),
    @(19, 44595, "Primera", 200, 15500, 16000, 15750, 875),
    @(20, 44984, "Especial", 160, 13000, 14000, 13500, 750),
    @(21, 44984, "Primera", 300, 11000, 12000, 11500, 639),
    @(22, 44630, "Especial", 300, 15000, 16000, 15500, 861),
    @(23, 44630, "Primera", 300, 12000, 13000, 12500, 694),
    @(24, 44630, "Segunda", 240, 9000, 10000, 9500, 528),
    @(25, 44295, "Especial", 200, 14500, 15000, 14750, 819),
    @(26, 44295, "Primera", 200, 12500, 13000, 12750, 708),
    @(27, 44295, "Segunda", 240, 10500, 11000, 10750, 597),
    @(28, 45001, "Especial", 400, 12000, 13000, 12500, 694),
    @(29, 45001, "Primera", 300, 10000, 11000, 10500, 583),
    @(30, 45009, "Especial", 300, 13000, 14000, 13500, 750),
    @(31, 45009, "Primera", 200, 11000, 12000, 11500, 639),
    @(32, 44294, "Especial", 200, 14500, 15000, 14750, 819),
    @(33, 44294, "Primera", 240, 12500, 13000, 12750, 708),
    @(34, 44294, "Segunda", 240, 10500, 11000, 10750, 597),
    @(35, 44637, "Especial", 200, 14000, 15000, 14500, 806),
    @(36, 44637, "Primera", 240, 10000, 11000, 10500, 583),
    @(37, 44987, "Especial", 160, 15000, 16000, 15500, 861),
    @(38, 44987, "Primera", 240, 12000, 13000, 12500, 694),
    @(39, 45015, "Especial", 200, 13000, 14000, 13500, 750),
    @(40, 45015, "Primera", 160, 11000, 12000, 11500, 639),
    @(41, 45005, "Especial", 300, 13000, 14000, 13500, 750),
    @(42, 45005, "Primera", 200, 11000, 12000, 11500, 639),
    @(43, 44603, "Especial", 240, 14500, 15000, 14750, 819),
    @(44, 44610, "Primera", 200, 13000, 14000, 13500, 750),
    @(45, 44610, "Segunda", 200, 11000, 12000, 11500, 639),
    @(46, 44636, "Especial", 240, 14000, 15000, 14500, 806),
    @(47, 44636, "Primera", 200, 10000, 11000, 10500, 583),
    @(48, 44606, "Primera", 240, 11500, 12000, 11750, 653),
    @(49, 44606, "Segunda", 240, 9500, 10000, 9750, 542),
    @(50, 44988, "Especial", 100, 14000, 15000, 14500, 806),
    @(51, 44988, "Primera", 200, 11000, 12000, 11500, 639)
)

# Row 51 is brand new: seed the columns that are constant for every data
# row in this sheet (A, B, C, E, F, G, H, I, J, K, Q, R, T) before filling
# in the per-row values below. Also copy the date number format used by
# column D onto the new row so the date keeps displaying correctly.
$ws.Cells.Item(51, 4).NumberFormat = $ws.Cells.Item(50, 4).NumberFormat
$ws.Cells.Item(51, 1).Value = $ws.Cells.Item(50, 1).Value2
$ws.Cells.Item(51, 2).Value = $ws.Cells.Item(50, 2).Value2
$ws.Cells.Item(51, 3).Value = $ws.Cells.Item(50, 3).Value2
$ws.Cells.Item(51, 5).Value = $ws.Cells.Item(50, 5).Value2
$ws.Cells.Item(51, 6).Value = $ws.Cells.Item(50, 6).Value2
$ws.Cells.Item(51, 7).Value = $ws.Cells.Item(50, 7).Value2
$ws.Cells.Item(51, 8).Value = $ws.Cells.Item(50, 8).Value2
$ws.Cells.Item(51, 9).Value = $ws.Cells.Item(50, 9).Value2
$ws.Cells.Item(51, 10).Value = $ws.Cells.Item(50, 10).Value2
$ws.Cells.Item(51, 11).Value = $ws.Cells.Item(50, 11).Value2
$ws.Cells.Item(51, 17).Value = $ws.Cells.Item(50, 17).Value2
$ws.Cells.Item(51, 18).Value = $ws.Cells.Item(50, 18).Value2
$ws.Cells.Item(51, 20).Value = $ws.Cells.Item(50, 20).Value2

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D - Fecha
    $ws.Cells.Item($r, 12).Value = $row[2]  # L - Calidad
    $ws.Cells.Item($r, 13).Value = $row[3]  # M - Volumen
    $ws.Cells.Item($r, 14).Value = $row[4]  # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $row[5]  # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $row[6]  # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $row[7]  # S - Precio $/Kg
}
